# Finished 3cm curvature analysis for 2015
# Populate AREA/MIN/MAX/RANGE/MEAN/STD/SUM columns (D:J) for rows 16-24
# on all three sheets, then update the saved view state (selection +
# active tab) to match where the author left off.

$wb = $excel.ActiveWorkbook

$wsCurv = $wb.Worksheets.Item("Curvature")
$wsPlan = $wb.Worksheets.Item("Plan_Curve")
$wsProf = $wb.Worksheets.Item("Profile_Curve")

# --- Curvature ---
$wsCurv.Range("D16").Value = 22096
$wsCurv.Range("E16").Value = 19.886399999999998
$wsCurv.Range("F16").Value = -113057.5625
$wsCurv.Range("G16").Value = 94997.0625
$wsCurv.Range("H16").Value = 208054.625
$wsCurv.Range("I16").Value = 20.760270999999999
$wsCurv.Range("J16").Value = 9277.4813969999996
$wsCurv.Range("D17").Value = 22314
$wsCurv.Range("E17").Value = 20.082599999999999
$wsCurv.Range("F17").Value = -79605
$wsCurv.Range("G17").Value = 50441.476562999997
$wsCurv.Range("H17").Value = 130046.476563
$wsCurv.Range("I17").Value = 4.2869099999999998
$wsCurv.Range("J17").Value = 6726.4234059999999
$wsCurv.Range("D18").Value = 16131
$wsCurv.Range("E18").Value = 14.517899999999999
$wsCurv.Range("F18").Value = -48516.960937999997
$wsCurv.Range("G18").Value = 44731.777344000002
$wsCurv.Range("H18").Value = 93248.738280999998
$wsCurv.Range("I18").Value = -12.991279
$wsCurv.Range("J18").Value = 6157.7164089999997
$wsCurv.Range("D19").Value = 19415
$wsCurv.Range("E19").Value = 17.473500000000001
$wsCurv.Range("F19").Value = -92679.898438000004
$wsCurv.Range("G19").Value = 93582.601563000004
$wsCurv.Range("H19").Value = 186262.5
$wsCurv.Range("I19").Value = -1.1193040000000001
$wsCurv.Range("J19").Value = 9120.9722089999996
$wsCurv.Range("D20").Value = 18835
$wsCurv.Range("E20").Value = 16.951499999999999
$wsCurv.Range("F20").Value = -78936.125
$wsCurv.Range("G20").Value = 119264.101563
$wsCurv.Range("H20").Value = 198200.226563
$wsCurv.Range("I20").Value = 19.950292999999999
$wsCurv.Range("J20").Value = 8421.8076760000004
$wsCurv.Range("D21").Value = 21400
$wsCurv.Range("E21").Value = 19.260000000000002
$wsCurv.Range("F21").Value = -143074.25
$wsCurv.Range("G21").Value = 105557.46875
$wsCurv.Range("H21").Value = 248631.71875
$wsCurv.Range("I21").Value = 14.081272
$wsCurv.Range("J21").Value = 9108.5376219999998
$wsCurv.Range("D22").Value = 18613
$wsCurv.Range("E22").Value = 16.7517
$wsCurv.Range("F22").Value = -64001.746094000002
$wsCurv.Range("G22").Value = 58183.300780999998
$wsCurv.Range("H22").Value = 122185.046875
$wsCurv.Range("I22").Value = 22.793627999999998
$wsCurv.Range("J22").Value = 6006.4417119999998
$wsCurv.Range("D23").Value = 19843
$wsCurv.Range("E23").Value = 17.858699999999999
$wsCurv.Range("F23").Value = -94722.539063000004
$wsCurv.Range("G23").Value = 82953.867188000004
$wsCurv.Range("H23").Value = 177676.40625
$wsCurv.Range("I23").Value = 19.693408000000002
$wsCurv.Range("J23").Value = 8616.8592590000007
$wsCurv.Range("D24").Value = 18845
$wsCurv.Range("E24").Value = 16.9605
$wsCurv.Range("F24").Value = -177645.546875
$wsCurv.Range("G24").Value = 120519.875
$wsCurv.Range("H24").Value = 298165.421875
$wsCurv.Range("I24").Value = 32.655202000000003
$wsCurv.Range("J24").Value = 12259.073715

# --- Plan_Curve ---
$wsPlan.Range("D16").Value = 22096
$wsPlan.Range("E16").Value = 19.886399999999998
$wsPlan.Range("F16").Value = -77939.507813000004
$wsPlan.Range("G16").Value = 63492.152344000002
$wsPlan.Range("H16").Value = 141431.660156
$wsPlan.Range("I16").Value = 87.873146000000006
$wsPlan.Range("J16").Value = 4347.5604960000001
$wsPlan.Range("D17").Value = 22314
$wsPlan.Range("E17").Value = 20.082599999999999
$wsPlan.Range("F17").Value = -37430.941405999998
$wsPlan.Range("G17").Value = 26780.464843999998
$wsPlan.Range("H17").Value = 64211.40625
$wsPlan.Range("I17").Value = 124.46781799999999
$wsPlan.Range("J17").Value = 3261.4611169999998
$wsPlan.Range("D18").Value = 16131
$wsPlan.Range("E18").Value = 14.517899999999999
$wsPlan.Range("F18").Value = -35157.300780999998
$wsPlan.Range("G18").Value = 26354.185547000001
$wsPlan.Range("H18").Value = 61511.486327999999
$wsPlan.Range("I18").Value = 132.71728300000001
$wsPlan.Range("J18").Value = 2951.3188719999998
$wsPlan.Range("D19").Value = 19415
$wsPlan.Range("E19").Value = 17.473500000000001
$wsPlan.Range("F19").Value = -53631.550780999998
$wsPlan.Range("G19").Value = 51617.335937999997
$wsPlan.Range("H19").Value = 105248.886719
$wsPlan.Range("I19").Value = 107.79254299999999
$wsPlan.Range("J19").Value = 4334.7871660000001
$wsPlan.Range("D20").Value = 18835
$wsPlan.Range("E20").Value = 16.951499999999999
$wsPlan.Range("F20").Value = -38377.035155999998
$wsPlan.Range("G20").Value = 57585.492187999997
$wsPlan.Range("H20").Value = 95962.527344000002
$wsPlan.Range("I20").Value = 113.21618100000001
$wsPlan.Range("J20").Value = 4045.466488
$wsPlan.Range("D21").Value = 21400
$wsPlan.Range("E21").Value = 19.260000000000002
$wsPlan.Range("F21").Value = -76341.359375
$wsPlan.Range("G21").Value = 62465.328125
$wsPlan.Range("H21").Value = 138806.6875
$wsPlan.Range("I21").Value = 122.871897
$wsPlan.Range("J21").Value = 4395.6796850000001
$wsPlan.Range("D22").Value = 18613
$wsPlan.Range("E22").Value = 16.7517
$wsPlan.Range("F22").Value = -31937.058593999998
$wsPlan.Range("G22").Value = 30746.246093999998
$wsPlan.Range("H22").Value = 62683.304687999997
$wsPlan.Range("I22").Value = 99.529325999999998
$wsPlan.Range("J22").Value = 2766.8014210000001
$wsPlan.Range("D23").Value = 19843
$wsPlan.Range("E23").Value = 17.858699999999999
$wsPlan.Range("F23").Value = -54001.914062999997
$wsPlan.Range("G23").Value = 59129.773437999997
$wsPlan.Range("H23").Value = 113131.6875
$wsPlan.Range("I23").Value = 136.55813599999999
$wsPlan.Range("J23").Value = 4023.072885
$wsPlan.Range("D24").Value = 18845
$wsPlan.Range("E24").Value = 16.9605
$wsPlan.Range("F24").Value = -111126.960938
$wsPlan.Range("G24").Value = 91346.828125
$wsPlan.Range("H24").Value = 202473.789063
$wsPlan.Range("I24").Value = 46.561487
$wsPlan.Range("J24").Value = 5972.5656269999999

# --- Profile_Curve ---
$wsProf.Range("D16").Value = 22096
$wsProf.Range("E16").Value = 19.886399999999998
$wsProf.Range("F16").Value = -61546.707030999998
$wsProf.Range("G16").Value = 62662.535155999998
$wsProf.Range("H16").Value = 124209.242188
$wsProf.Range("I16").Value = 67.112870999999998
$wsProf.Range("J16").Value = 6063.3135350000002
$wsProf.Range("D17").Value = 22314
$wsProf.Range("E17").Value = 20.082599999999999
$wsProf.Range("F17").Value = -33802.496094000002
$wsProf.Range("G17").Value = 42174.054687999997
$wsProf.Range("H17").Value = 75976.550780999998
$wsProf.Range("I17").Value = 120.180908
$wsProf.Range("J17").Value = 4283.4623220000003
$wsProf.Range("D18").Value = 16131
$wsProf.Range("E18").Value = 14.517899999999999
$wsProf.Range("F18").Value = -24172.005859000001
$wsProf.Range("G18").Value = 28797.314452999999
$wsProf.Range("H18").Value = 52969.320312999997
$wsProf.Range("I18").Value = 145.708564
$wsProf.Range("J18").Value = 3973.0706340000002
$wsProf.Range("D19").Value = 19415
$wsProf.Range("E19").Value = 17.473500000000001
$wsProf.Range("F19").Value = -52245.726562999997
$wsProf.Range("G19").Value = 48715.410155999998
$wsProf.Range("H19").Value = 100961.136719
$wsProf.Range("I19").Value = 108.91184699999999
$wsProf.Range("J19").Value = 5881.8981130000002
$wsProf.Range("D20").Value = 18835
$wsProf.Range("E20").Value = 16.951499999999999
$wsProf.Range("F20").Value = -61678.605469000002
$wsProf.Range("G20").Value = 40559.089844000002
$wsProf.Range("H20").Value = 102237.695313
$wsProf.Range("I20").Value = 93.265889000000001
$wsProf.Range("J20").Value = 5405.6438870000002
$wsProf.Range("D21").Value = 21400
$wsProf.Range("E21").Value = 19.260000000000002
$wsProf.Range("F21").Value = -60859.414062999997
$wsProf.Range("G21").Value = 66732.890625
$wsProf.Range("H21").Value = 127592.304688
$wsProf.Range("I21").Value = 108.79062500000001
$wsProf.Range("J21").Value = 5818.3414240000002
$wsProf.Range("D22").Value = 18613
$wsProf.Range("E22").Value = 16.7517
$wsProf.Range("F22").Value = -39631.742187999997
$wsProf.Range("G22").Value = 48100.265625
$wsProf.Range("H22").Value = 87732.007813000004
$wsProf.Range("I22").Value = 76.735697000000002
$wsProf.Range("J22").Value = 4036.4183560000001
$wsProf.Range("D23").Value = 19843
$wsProf.Range("E23").Value = 17.858699999999999
$wsProf.Range("F23").Value = -47754.503905999998
$wsProf.Range("G23").Value = 47029.464844000002
$wsProf.Range("H23").Value = 94783.96875
$wsProf.Range("I23").Value = 116.86472500000001
$wsProf.Range("J23").Value = 5743.853521
$wsProf.Range("D24").Value = 18845
$wsProf.Range("E24").Value = 16.9605
$wsProf.Range("F24").Value = -59257.761719000002
$wsProf.Range("G24").Value = 84534.875
$wsProf.Range("H24").Value = 143792.636719
$wsProf.Range("I24").Value = 13.906283999999999
$wsProf.Range("J24").Value = 7789.2204760000004

# Update the saved selection on every sheet, and leave Profile_Curve as
# the active (visible) tab, matching the workbook's activeTab state.
$wsCurv.Range("A25").Select()
$wsPlan.Range("A25").Select()
$wsProf.Range("A25").Select()
$wsProf.Activate()
